$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "562.74"
# or "138.30" are not auto-converted to numbers by Excel's smart parsing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.140.73"
$ws.Range("E2").Value = "  +3.61%  "

$ws.Range("D3").Value = "2.989.81"
$ws.Range("E3").Value = "  +3.36%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "562.74"
$ws.Range("E5").Value = "  +3.01%  "

$ws.Range("D6").Value = "138.30"
$ws.Range("E6").Value = "  +11.26%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +3.86%  "

$ws.Range("D9").Value = "2.979.19"
$ws.Range("E9").Value = "  +2.99%  "

$ws.Range("E10").Value = "  +8.34%  "

$ws.Range("D11").Value = "5.08"
$ws.Range("E11").Value = "  +9.36%  "

$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  +5.29%  "

$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +9.93%  "

$ws.Range("D14").Value = "33.62"
$ws.Range("E14").Value = "  +5.07%  "

$ws.Range("E15").Value = "  +2.75%  "

$ws.Range("D16").Value = "3.482.64"
$ws.Range("E16").Value = "  +3.39%  "

$ws.Range("D17").Value = "7.03"
$ws.Range("E17").Value = "  +8.77%  "

$ws.Range("D18").Value = "2.986.17"
$ws.Range("E18").Value = "  +3.45%  "

$ws.Range("D19").Value = "59.100.98"
$ws.Range("E19").Value = "  +3.58%  "

$ws.Range("D20").Value = "426.84"
$ws.Range("E20").Value = "  +6.21%  "

$ws.Range("D21").Value = "13.53"
$ws.Range("E21").Value = "  +6.35%  "

$ws.Range("D22").Value = "0.713"
$ws.Range("E22").Value = "  +7.18%  "

$ws.Range("D23").Value = "13.43"
$ws.Range("E23").Value = "  +6.76%  "

$ws.Range("E24").Value = "  +4.35%  "

$ws.Range("D25").Value = "80.59"
$ws.Range("E25").Value = "  +4.22%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "2.13"
$ws.Range("E28").Value = "  +12.08%  "

$ws.Range("E29").Value = "  +3.70%  "

$ws.Range("E30").Value = "  +7.56%  "

$ws.Range("E31").Value = "  +4.63%  "

$ws.Range("D32").Value = "6.11"
$ws.Range("E32").Value = "  +3.42%  "

$ws.Range("D33").Value = "0.0984"
$ws.Range("E33").Value = "  +0.74%  "

$ws.Range("D34").Value = "0.0₃0770"
$ws.Range("E34").Value = "  +23.76%  "

$ws.Range("D35").Value = "0.990"
$ws.Range("E35").Value = "  +8.66%  "

$ws.Range("E36").Value = "  +7.58%  "

$ws.Range("E37").Value = "  +3.65%  "

$ws.Range("D38").Value = "49.08"
$ws.Range("E38").Value = "  +2.13%  "

$ws.Range("D39").Value = "8.64"
$ws.Range("E39").Value = "  +5.51%  "

$ws.Range("E40").Value = "  +13.57%  "

$ws.Range("D41").Value = "398.48"
$ws.Range("E41").Value = "  +10.62%  "

$ws.Range("E42").Value = "  +4.43%  "

$ws.Range("D43").Value = "2.747.55"
$ws.Range("E43").Value = "  +5.35%  "

$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("E45").Value = "  +10.70%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "122.80"
$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("E48").Value = "  +2.64%  "

$ws.Range("E49").Value = "  +4.11%  "

$ws.Range("D50").Value = "32.58"
$ws.Range("E50").Value = "  +20.74%  "

$ws.Range("D51").Value = "23.42"
$ws.Range("E51").Value = "  +4.55%  "

# Restore the original (default/"Normal") style on column D so only the
# cell text content changes, matching the source diff (no formatting deltas).
$ws.Range("D2:D51").Style = "Normal"
